$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Swap the OLE_LINK103 / OLE_LINK104 bookmark names (they wrap the
#    same span of text, so this just flips which bookmark object owns
#    which name). Bookmark.Name is read-only in the Word OM, so the
#    way to "rename" is to re-Add a bookmark with the desired name at
#    the target range -- Word moves/redefines it rather than
#    duplicating when the name already exists.
# ------------------------------------------------------------------
$bm103 = $d.Bookmarks.Item("OLE_LINK103")
$bm104 = $d.Bookmarks.Item("OLE_LINK104")
$r103 = $bm103.Range
$r104 = $bm104.Range
$d.Bookmarks.Add("OLE_LINK103", $r104)
$d.Bookmarks.Add("OLE_LINK104", $r103)

# ------------------------------------------------------------------
# 2. "期望：..." expectation sentence rewrite.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "期望：总转化率、留存率、净转化率都减小，都具有统计显著性与实际显著性",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "期望：总转化率、留存率增加，都具有统计显著性与实际显著性，净转化率，变化不大。", 1) | Out-Null

# ------------------------------------------------------------------
# 3. "通过在线计算器计算出：" -> "通过在线计算器与所给数据计算出："
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "通过在线计算器计算出：",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "通过在线计算器与所给数据计算出：", 1) | Out-Null

# ------------------------------------------------------------------
# 4. Standalone number 4737818 -> 645875 (the first/standalone run,
#    which occurs earlier in the document than the sentence in step 5).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "4737818",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "645875", 1) | Out-Null

# ------------------------------------------------------------------
# 5. "选择总转化率网页计算天数：..." sentence rewrite (still has the old
#    4737818 text embedded, so match on the full old sentence).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "选择总转化率网页计算天数：4737818/40000 =118.45 大概需要4个月的时间，持续时间太长",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "选择总转化率网页计算天数：645875/40000 =16.15 大概需要17天时间", 1) | Out-Null

# ------------------------------------------------------------------
# 6. "选择净转化率网页计算天数：..." sentence gets an extra clause inserted.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "选择净转化率网页计算天数：685325/40000 = 17.13 时间合适，如果四舍五入为17天，所需网页数不够，因此进位，选择持续时间18天",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "选择净转化率网页计算天数：685325/40000 = 17.13 时间合适，如果四舍五入为17天，所需网页数不够，不能计算净转化率，因此进位，选择持续时间18天", 1) | Out-Null

# ------------------------------------------------------------------
# 7. "转移单位：完成登录的用户Id数量" -> "转移单位：每个完成登录的用户Id"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "转移单位：完成登录的用户Id数量",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "转移单位：每个完成登录的用户Id", 1) | Out-Null

# ------------------------------------------------------------------
# 8. Move the hidden "_GoBack" bookmark (tracks the last-edit location)
#    from its old spot (after the "不具有统计显著性" paragraph) to the
#    blank paragraph two paragraphs after the "转移单位" line we just
#    edited -- re-Adding it at the new range moves it automatically.
# ------------------------------------------------------------------
$locator = $d.Content
$locator.Find.Execute(
    "转移单位：每个完成登录的用户Id",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$wdParagraph = 4
$targetPara = $locator.Next($wdParagraph, 2)
$d.Bookmarks.Add("_GoBack", $targetPara)
